$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A width: 19.28515625 -> 24.0 ---
$ws.Columns.Item(1).ColumnWidth = 23.166666666666668

# Seed rows 20 and 21 with the same formatting as row 19 (the existing
# "annotation user" row immediately above) before overwriting the values,
# so the new rows pick up style s="6" / s="7" like their neighbours.
$ws.Range("A19:G19").Copy()
$ws.Range("A20:G20").Select()
$excel.ActiveSheet.Paste()
$ws.Range("A19:G19").Copy()
$ws.Range("A21:G21").Select()
$excel.ActiveSheet.Paste()

# --- Row 20: new user "shareAnnotationUser1" ---
$ws.Range("A20").Value = "shareAnnotationUser1"
$ws.Range("B20").Value = "Password1"
$ws.Range("C20").Value = ""
$ws.Range("D20").Value = ""
$ws.Range("E20").Value = "ANZ annotation user"
$ws.Range("F20").Value = ""
$ws.Range("G20").Value = "shareannotationuser1@mailinator.com"
$ws.Hyperlinks.Add($ws.Range("G20"), "mailto:shareannotationuser1@mailinator.com")

# --- Row 21: new user "myShareAnnotationUser" ---
$ws.Range("A21").Value = "myShareAnnotationUser"
$ws.Range("B21").Value = "Password1"
$ws.Range("C21").Value = ""
$ws.Range("D21").Value = ""
$ws.Range("E21").Value = "ANZ annotation user"
$ws.Range("F21").Value = ""
$ws.Range("G21").Value = "myShareAnnotationUser@mailinator.com"
$ws.Hyperlinks.Add($ws.Range("G21"), "mailto:myShareAnnotationUser@mailinator.com")

# --- Selection moves to J24 ---
$ws.Range("J24").Select()
